$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("G20").Value = "..."
$ws.Range("H20").Value = "etc"
$ws.Range("I20").Value = "..."
$ws.Range("G21").Value = "..."
$ws.Range("H21").Value = "etc"
$ws.Range("I21").Value = "..."
$ws.Range("G20:I21").HorizontalAlignment = -4108
